# The deck ships two embedded themes:
#   ppt/theme/theme1.xml  - "Office Theme" (clrScheme name "Office")
#   ppt/theme/theme2.xml  - "Integral" / "Red Violet" (the one actually
#                            applied to the slide master / whole deck)
#
# The authored edit swaps the two theme parts so the deck's live design
# becomes the plain "Office Theme" colour scheme instead of "Integral".
# Re-create that swap through the PowerPoint object model by pushing the
# Office theme's twelve scheme colours onto the presentation's active
# theme (Slide Master's Theme.ThemeColorScheme) - font/format schemes in
# both theme parts are already identical, so the colour swap is the only
# visible change that needs to be applied.

function Get-BGRLong($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# Target palette = the "Office Theme" colours currently sitting unused in
# theme1.xml, in MsoThemeColorSchemeIndex order (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink).
$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$scheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $scheme.Colors($i).RGB = Get-BGRLong $officeColors[$i - 1]
}
